$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "70.1% [68.6%;71.8%]"
$ws.Range("F2").Value = "70.1% [68.6%;71.8%]"
$ws.Range("G2").Value = "70.1% [68.6%;71.8%]"
$ws.Range("I2").Value = "11.1% [9.9%;12.2%]"
$ws.Range("J2").Value = "95.0% [90.0%;99.0%]"
$ws.Range("E3").Value = "71.9% [70.3%;73.5%]"
$ws.Range("F3").Value = "71.9% [70.3%;73.6%]"
$ws.Range("G3").Value = "71.9% [70.3%;73.5%]"
$ws.Range("I3").Value = "5.2% [4.3%;6.1%]"
$ws.Range("J3").Value = "72.0% [63.0%;80.0%]"
$ws.Range("E4").Value = "85.5% [83.0%;87.9%]"
$ws.Range("F4").Value = "64.0% [61.9%;66.1%]"
$ws.Range("I4").Value = "10.3% [8.7%;12.1%]"
$ws.Range("J4").Value = "61.0% [51.0%;70.0%]"
$ws.Range("E5").Value = "82.6% [80.5%;84.7%]"
$ws.Range("F5").Value = "82.6% [80.5%;84.7%]"
$ws.Range("G5").Value = "82.6% [80.5%;84.6%]"
$ws.Range("I5").Value = "0.8% [0.6%;1.1%]"
$ws.Range("J5").Value = "49.0% [39.0%;59.0%]"
$ws.Range("E6").Value = "71.6% [70.1%;73.1%]"
$ws.Range("F6").Value = "71.6% [70.1%;73.1%]"
$ws.Range("G6").Value = "71.6% [70.1%;73.1%]"
$ws.Range("I6").Value = "0.1% [0.1%;0.2%]"
$ws.Range("J6").Value = "13.0% [7.0%;20.0%]"
$ws.Range("E7").Value = "89.7% [87.6%;91.7%]"
$ws.Range("F7").Value = "71.3% [69.7%;73.0%]"
$ws.Range("G7").Value = "66.3% [65.2%;67.4%]"
$ws.Range("H7").Value = "5.0% [5.0%;5.1%]"
$ws.Range("I7").Value = "7.1% [6.4%;7.9%]"
$ws.Range("J7").Value = "63.0% [53.0%;72.0%]"
$ws.Range("E8").Value = "98.2% [97.0%;99.2%]"
$ws.Range("F8").Value = "98.2% [97.0%;99.2%]"
$ws.Range("G8").Value = "98.2% [97.0%;99.2%]"
$ws.Range("H8").Value = "11.8% [9.2%;14.6%]"
$ws.Range("I8").Value = "9.4% [7.4%;11.6%]"
$ws.Range("E9").Value = "97.2% [95.6%;98.6%]"
$ws.Range("F9").Value = "97.2% [95.6%;98.6%]"
$ws.Range("G9").Value = "97.2% [95.6%;98.6%]"
$ws.Range("H9").Value = "8.0% [5.8%;10.2%]"
$ws.Range("E10").Value = "95.6% [93.8%;97.2%]"
$ws.Range("F10").Value = "93.2% [91.0%;95.2%]"
$ws.Range("G10").Value = "88.0% [85.4%;90.6%]"
$ws.Range("H10").Value = "5.0% [3.2%;7.2%]"
$ws.Range("I10").Value = "2.2% [1.0%;3.5%]"
$ws.Range("J10").Value = "7.0% [2.0%;12.0%]"
$ws.Range("F11").Value = "92.2% [90.0%;94.4%]"
$ws.Range("G11").Value = "87.8% [85.2%;90.2%]"
$ws.Range("H11").Value = "3.8% [2.2%;5.6%]"
$ws.Range("I11").Value = "2.1% [0.9%;3.4%]"
$ws.Range("J11").Value = "6.0% [2.0%;11.0%]"
$ws.Range("H12").Value = "17.2% [13.9%;20.7%]"
$ws.Range("I12").Value = "63.0% [53.0%;72.0%]"
$ws.Range("J12").Value = "63.0% [53.0%;72.0%]"
$ws.Range("H13").Value = "16.5% [13.3%;19.8%]"
$ws.Range("I13").Value = "64.0% [54.0%;73.0%]"
$ws.Range("J13").Value = "64.0% [54.0%;73.0%]"
$ws.Range("H14").Value = "15.1% [12.6%;17.7%]"
$ws.Range("I14").Value = "14.0% [8.0%;21.0%]"
$ws.Range("J14").Value = "11.0% [5.0%;17.0%]"
$ws.Range("H15").Value = "13.8% [11.6%;16.2%]"
$ws.Range("I15").Value = "14.0% [8.0%;21.0%]"
$ws.Range("J15").Value = "13.0% [7.0%;20.0%]"
$ws.Range("E18").Value = "92.0% [86.0%;97.0%]"
$ws.Range("F18").Value = "82.0% [74.0%;89.0%]"
$ws.Range("G18").Value = "82.0% [74.0%;89.0%]"
$ws.Range("H18").Value = "5.5% [3.2%;7.8%]"
$ws.Range("I18").Value = "6.0% [2.5%;10.0%]"
$ws.Range("J18").Value = "7.0% [2.0%;12.0%]"
$ws.Range("E19").Value = "90.0% [84.0%;95.0%]"
$ws.Range("F19").Value = "80.0% [72.0%;88.0%]"
$ws.Range("G19").Value = "80.0% [72.0%;88.0%]"
$ws.Range("H19").Value = "6.0% [3.8%;8.5%]"
$ws.Range("I19").Value = "5.7% [2.0%;10.0%]"
$ws.Range("J19").Value = "7.0% [2.0%;12.0%]"
